$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 2233.6667  # H2: 1763.5454 -> 2233.6667
$ws.Cells.Item(2, 9).Value = 2950  # I2: 1977.6666 -> 2950
$ws.Cells.Item(2, 10).Value = 801  # J2: 800 -> 801
$ws.Cells.Item(2, 11).Value = 2950  # K2: 1977.6666 -> 2950
$ws.Cells.Item(2, 12).Value = 801  # L2: 800 -> 801
$ws.Cells.Item(2, 13).Value = -2837  # M2: -1864.6666 -> -2837
$ws.Cells.Item(2, 14).Value = -1027  # N2: -1026 -> -1027
# Row 9
$ws.Cells.Item(9, 8).Value = 275.66666  # H9: 246.75 -> 275.66666
$ws.Cells.Item(9, 10).Value = 338.5  # J9: 279 -> 338.5
$ws.Cells.Item(9, 12).Value = 338.5  # L9: 279 -> 338.5
$ws.Cells.Item(9, 14).Value = -676.5  # N9: -617 -> -676.5
# Row 21
$ws.Cells.Item(21, 8).Value = 8335820  # H21: 9095397 -> 8335820
$ws.Cells.Item(21, 9).Value = 16667607  # I21: 20002384 -> 16667607
$ws.Cells.Item(21, 10).Value = 4033.1667  # J21: 6241.3335 -> 4033.1667
$ws.Cells.Item(21, 11).Value = 16667607  # K21: 20002384 -> 16667607
$ws.Cells.Item(21, 12).Value = 4033.1667  # L21: 6241.3335 -> 4033.1667
$ws.Cells.Item(21, 13).Value = -16667139  # M21: -20001916 -> -16667139
$ws.Cells.Item(21, 14).Value = -4969.1667  # N21: -7177.3335 -> -4969.1667
# Row 23
$ws.Cells.Item(23, 8).Value = 8335820  # H23: 9095397 -> 8335820
$ws.Cells.Item(23, 9).Value = 16667607  # I23: 20002384 -> 16667607
$ws.Cells.Item(23, 10).Value = 4033.1667  # J23: 6241.3335 -> 4033.1667
$ws.Cells.Item(23, 11).Value = 16667607  # K23: 20002384 -> 16667607
$ws.Cells.Item(23, 12).Value = 4033.1667  # L23: 6241.3335 -> 4033.1667
$ws.Cells.Item(23, 13).Value = -16667373  # M23: -20002150 -> -16667373
$ws.Cells.Item(23, 14).Value = -4501.1667  # N23: -6709.3335 -> -4501.1667
# Row 48
$ws.Cells.Item(48, 8).Value = 11750  # H48: 13500 -> 11750
$ws.Cells.Item(48, 9).Value = 10999.5  # I48: 12000 -> 10999.5
$ws.Cells.Item(48, 10).Value = 12500.5  # J48: 15000 -> 12500.5
$ws.Cells.Item(48, 11).Value = 32998.5  # K48: 36000 -> 32998.5
$ws.Cells.Item(48, 12).Value = 37501.5  # L48: 45000 -> 37501.5
$ws.Cells.Item(48, 13).Value = -32706.5  # M48: -35708 -> -32706.5
$ws.Cells.Item(48, 14).Value = -38085.5  # N48: -45584 -> -38085.5
# Row 56
$ws.Cells.Item(56, 8).Value = 11750  # H56: 13500 -> 11750
$ws.Cells.Item(56, 9).Value = 10999.5  # I56: 12000 -> 10999.5
$ws.Cells.Item(56, 10).Value = 12500.5  # J56: 15000 -> 12500.5
$ws.Cells.Item(56, 11).Value = 32998.5  # K56: 36000 -> 32998.5
$ws.Cells.Item(56, 12).Value = 37501.5  # L56: 45000 -> 37501.5
$ws.Cells.Item(56, 13).Value = -32464.5  # M56: -35466 -> -32464.5
$ws.Cells.Item(56, 14).Value = -38569.5  # N56: -46068 -> -38569.5
# Row 74
$ws.Cells.Item(74, 8).Value = 4141.95  # H74: 4441.0557 -> 4141.95
$ws.Cells.Item(74, 9).Value = 2921.3572  # I74: 3166.5833 -> 2921.3572
$ws.Cells.Item(74, 11).Value = 2921.3572  # K74: 3166.5833 -> 2921.3572
$ws.Cells.Item(74, 13).Value = -1985.3572  # M74: -2230.5833 -> -1985.3572
# Row 77
$ws.Cells.Item(77, 8).Value = 4141.95  # H77: 4441.0557 -> 4141.95
$ws.Cells.Item(77, 9).Value = 2921.3572  # I77: 3166.5833 -> 2921.3572
$ws.Cells.Item(77, 11).Value = 14606.786  # K77: 15832.9165 -> 14606.786
$ws.Cells.Item(77, 13).Value = -9926.786  # M77: -11152.9165 -> -9926.786
# Row 112
$ws.Cells.Item(112, 8).Value = 3850.6667  # H112: 4019.0908 -> 3850.6667
$ws.Cells.Item(112, 10).Value = 3734.2222  # J112: 3951.25 -> 3734.2222
$ws.Cells.Item(112, 12).Value = 11202.6666  # L112: 11853.75 -> 11202.6666
$ws.Cells.Item(112, 14).Value = -13418.6666  # N112: -14069.75 -> -13418.6666
# Row 116
$ws.Cells.Item(116, 8).Value = 4587.5  # H116: 5450 -> 4587.5
$ws.Cells.Item(116, 9).Value = 2175  # I116: 2350 -> 2175
$ws.Cells.Item(116, 11).Value = 2175  # K116: 2350 -> 2175
$ws.Cells.Item(116, 13).Value = 1267  # M116: 1092 -> 1267
# Row 138
$ws.Cells.Item(138, 8).Value = 3328.4666  # H138: 3465.2593 -> 3328.4666
$ws.Cells.Item(138, 10).Value = 2892.6667  # J138: 3051.7334 -> 2892.6667
$ws.Cells.Item(138, 12).Value = 8678.000100000001  # L138: 9155.200199999999 -> 8678.000100000001
$ws.Cells.Item(138, 14).Value = -18958.0001  # N138: -19435.2002 -> -18958.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Cells.Item(10, 8).Value = 10060  # H10: 50010000 -> 10060
$ws.Cells.Item(10, 10).Value = 120  # J10: 100000000 -> 120
$ws.Cells.Item(10, 12).Value = 120  # L10: 100000000 -> 120
$ws.Cells.Item(10, 14).Value = -460  # N10: -100000340 -> -460
# Row 11
$ws.Cells.Item(11, 8).Value = 200  # H11: 975 -> 200
$ws.Cells.Item(11, 9).Value = 200  # I11: 975 -> 200
$ws.Cells.Item(11, 11).Value = 200  # K11: 975 -> 200
$ws.Cells.Item(11, 13).Value = -56  # M11: -831 -> -56
# Row 13
$ws.Cells.Item(13, 10).Value = 0  # J13: 200 -> 0
$ws.Cells.Item(13, 12).Value = 0  # L13: 200 -> 0
$ws.Cells.Item(13, 14).ClearContents()  # N13: -488 -> (blank)
# Row 61
$ws.Cells.Item(61, 8).Value = 6365.2666  # H61: 6462.7856 -> 6365.2666
$ws.Cells.Item(61, 10).Value = 5750  # J61: 6500 -> 5750
$ws.Cells.Item(61, 12).Value = 5750  # L61: 6500 -> 5750
$ws.Cells.Item(61, 14).Value = -6174  # N61: -6924 -> -6174
# Row 68
$ws.Cells.Item(68, 8).Value = 60000  # H68: 0 -> 60000
$ws.Cells.Item(68, 10).Value = 60000  # J68: 0 -> 60000
$ws.Cells.Item(68, 12).Value = 60000  # L68: 0 -> 60000
$ws.Cells.Item(68, 14).Value = -61622  # N68: None -> -61622
# Row 71
$ws.Cells.Item(71, 8).Value = 60000  # H71: 0 -> 60000
$ws.Cells.Item(71, 10).Value = 60000  # J71: 0 -> 60000
$ws.Cells.Item(71, 12).Value = 180000  # L71: 0 -> 180000
$ws.Cells.Item(71, 14).Value = -188112  # N71: None -> -188112
# Row 125
$ws.Cells.Item(125, 8).Value = 0  # H125: 59999 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 59999 -> 0
$ws.Cells.Item(125, 12).Value = 0  # L125: 59999 -> 0
$ws.Cells.Item(125, 14).ClearContents()  # N125: -69839 -> (blank)
# Row 136
$ws.Cells.Item(136, 8).Value = 6365.2666  # H136: 6462.7856 -> 6365.2666
$ws.Cells.Item(136, 10).Value = 5750  # J136: 6500 -> 5750
$ws.Cells.Item(136, 12).Value = 17250  # L136: 19500 -> 17250
$ws.Cells.Item(136, 14).Value = -22350  # N136: -24600 -> -22350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Cells.Item(33, 8).Value = 2833.3333  # H33: 9750 -> 2833.3333
$ws.Cells.Item(33, 9).Value = 2833.3333  # I33: 4500 -> 2833.3333
$ws.Cells.Item(33, 10).Value = 0  # J33: 11500 -> 0
$ws.Cells.Item(33, 11).Value = 2833.3333  # K33: 4500 -> 2833.3333
$ws.Cells.Item(33, 12).Value = 0  # L33: 11500 -> 0
$ws.Cells.Item(33, 13).Value = -2497.3333  # M33: -4164 -> -2497.3333
$ws.Cells.Item(33, 14).ClearContents()  # N33: -12172 -> (blank)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Cells.Item(4, 8).Value = 6785.857  # H4: 15051 -> 6785.857
$ws.Cells.Item(4, 10).Value = 6785.857  # J4: 15051 -> 6785.857
$ws.Cells.Item(4, 12).Value = 6785.857  # L4: 15051 -> 6785.857
$ws.Cells.Item(4, 14).Value = -7009.857  # N4: -15275 -> -7009.857
# Row 31
$ws.Cells.Item(31, 8).Value = 4921.3335  # H31: 5246.25 -> 4921.3335
$ws.Cells.Item(31, 9).Value = 3438.25  # I31: 3772.2856 -> 3438.25
$ws.Cells.Item(31, 10).Value = 7887.5  # J31: 7309.8 -> 7887.5
$ws.Cells.Item(31, 11).Value = 3438.25  # K31: 3772.2856 -> 3438.25
$ws.Cells.Item(31, 12).Value = 7887.5  # L31: 7309.8 -> 7887.5
$ws.Cells.Item(31, 13).Value = -3143.25  # M31: -3477.2856 -> -3143.25
$ws.Cells.Item(31, 14).Value = -8477.5  # N31: -7899.8 -> -8477.5
# Row 34
$ws.Cells.Item(34, 8).Value = 4921.3335  # H34: 5246.25 -> 4921.3335
$ws.Cells.Item(34, 9).Value = 3438.25  # I34: 3772.2856 -> 3438.25
$ws.Cells.Item(34, 10).Value = 7887.5  # J34: 7309.8 -> 7887.5
$ws.Cells.Item(34, 11).Value = 3438.25  # K34: 3772.2856 -> 3438.25
$ws.Cells.Item(34, 12).Value = 7887.5  # L34: 7309.8 -> 7887.5
$ws.Cells.Item(34, 13).Value = -3236.25  # M34: -3570.2856 -> -3236.25
$ws.Cells.Item(34, 14).Value = -8291.5  # N34: -7713.8 -> -8291.5
# Row 69
$ws.Cells.Item(69, 8).Value = 20833.166  # H69: 23224.75 -> 20833.166
$ws.Cells.Item(69, 9).Value = 18999.8  # I69: 22666.334 -> 18999.8
$ws.Cells.Item(69, 10).Value = 30000  # J69: 24900 -> 30000
$ws.Cells.Item(69, 11).Value = 18999.8  # K69: 22666.334 -> 18999.8
$ws.Cells.Item(69, 12).Value = 30000  # L69: 24900 -> 30000
$ws.Cells.Item(69, 13).Value = -18250.8  # M69: -21917.334 -> -18250.8
$ws.Cells.Item(69, 14).Value = -31498  # N69: -26398 -> -31498
# Row 72
$ws.Cells.Item(72, 8).Value = 20833.166  # H72: 23224.75 -> 20833.166
$ws.Cells.Item(72, 9).Value = 18999.8  # I72: 22666.334 -> 18999.8
$ws.Cells.Item(72, 10).Value = 30000  # J72: 24900 -> 30000
$ws.Cells.Item(72, 11).Value = 56999.39999999999  # K72: 67999.00199999999 -> 56999.39999999999
$ws.Cells.Item(72, 12).Value = 90000  # L72: 74700 -> 90000
$ws.Cells.Item(72, 13).Value = -53255.39999999999  # M72: -64255.00199999999 -> -53255.39999999999
$ws.Cells.Item(72, 14).Value = -97488  # N72: -82188 -> -97488
# Row 134
$ws.Cells.Item(134, 8).Value = 1770.5454  # H134: 1585.1333 -> 1770.5454
$ws.Cells.Item(134, 9).Value = 1358.5555  # I134: 1271.3846 -> 1358.5555
$ws.Cells.Item(134, 11).Value = 4075.6665  # K134: 3814.1538 -> 4075.6665
$ws.Cells.Item(134, 13).Value = -1540.6665  # M134: -1279.1538 -> -1540.6665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Cells.Item(32, 8).Value = 1250.6666  # H32: 1666.6666 -> 1250.6666
$ws.Cells.Item(32, 9).Value = 1300.8  # I32: 1666.6666 -> 1300.8
$ws.Cells.Item(32, 10).Value = 1000  # J32: 0 -> 1000
$ws.Cells.Item(32, 11).Value = 3902.4  # K32: 4999.9998 -> 3902.4
$ws.Cells.Item(32, 12).Value = 3000  # L32: 0 -> 3000
$ws.Cells.Item(32, 13).Value = -3619.4  # M32: -4716.9998 -> -3619.4
$ws.Cells.Item(32, 14).Value = -3566  # N32: None -> -3566
# Row 34
$ws.Cells.Item(34, 8).Value = 1252.9642  # H34: 1258.6296 -> 1252.9642
$ws.Cells.Item(34, 10).Value = 3229.4  # J34: 3761.75 -> 3229.4
$ws.Cells.Item(34, 12).Value = 9688.200000000001  # L34: 11285.25 -> 9688.200000000001
$ws.Cells.Item(34, 14).Value = -9856.200000000001  # N34: -11453.25 -> -9856.200000000001
# Row 38
$ws.Cells.Item(38, 8).Value = 107.8421  # H38: 112.63158 -> 107.8421
$ws.Cells.Item(38, 9).Value = 54.25  # I38: 57.266666 -> 54.25
$ws.Cells.Item(38, 10).Value = 393.66666  # J38: 320.25 -> 393.66666
$ws.Cells.Item(38, 11).Value = 162.75  # K38: 171.799998 -> 162.75
$ws.Cells.Item(38, 12).Value = 1180.99998  # L38: 960.75 -> 1180.99998
$ws.Cells.Item(38, 13).Value = 184.25  # M38: 175.200002 -> 184.25
$ws.Cells.Item(38, 14).Value = -1874.99998  # N38: -1654.75 -> -1874.99998
# Row 54
$ws.Cells.Item(54, 8).Value = 444  # H54: 440 -> 444
$ws.Cells.Item(54, 10).Value = 444  # J54: 440 -> 444
$ws.Cells.Item(54, 12).Value = 1332  # L54: 1320 -> 1332
$ws.Cells.Item(54, 14).Value = -2450  # N54: -2438 -> -2450
# Row 121
$ws.Cells.Item(121, 8).Value = 873.36365  # H121: 939.5 -> 873.36365
$ws.Cells.Item(121, 10).Value = 1089  # J121: 1214.2858 -> 1089
$ws.Cells.Item(121, 12).Value = 3267  # L121: 3642.8574 -> 3267
$ws.Cells.Item(121, 14).Value = -5887  # N121: -6262.857400000001 -> -5887
# Row 129
$ws.Cells.Item(129, 8).Value = 1424.6666  # H129: 1579.6 -> 1424.6666
$ws.Cells.Item(129, 9).Value = 599.6667  # I129: 633 -> 599.6667
$ws.Cells.Item(129, 10).Value = 2249.6667  # J129: 2999.5 -> 2249.6667
$ws.Cells.Item(129, 11).Value = 1799.0001  # K129: 1899 -> 1799.0001
$ws.Cells.Item(129, 12).Value = 6749.000100000001  # L129: 8998.5 -> 6749.000100000001
$ws.Cells.Item(129, 13).Value = 3200.9999  # M129: 3101 -> 3200.9999
$ws.Cells.Item(129, 14).Value = -16749.0001  # N129: -18998.5 -> -16749.0001
# Row 131
$ws.Cells.Item(131, 8).Value = 2388.7817  # H131: 2401.8147 -> 2388.7817
$ws.Cells.Item(131, 10).Value = 2528.9387  # J131: 2546.5208 -> 2528.9387
$ws.Cells.Item(131, 12).Value = 7586.8161  # L131: 7639.562399999999 -> 7586.8161
$ws.Cells.Item(131, 14).Value = -17666.8161  # N131: -17719.5624 -> -17666.8161

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 2281.5  # H5: 2238 -> 2281.5
$ws.Cells.Item(5, 10).Value = 2499.5  # J5: 2500 -> 2499.5
$ws.Cells.Item(5, 12).Value = 2499.5  # L5: 2500 -> 2499.5
$ws.Cells.Item(5, 14).Value = -2723.5  # N5: -2724 -> -2723.5
# Row 59
$ws.Cells.Item(59, 8).Value = 5000  # H59: 3000 -> 5000
$ws.Cells.Item(59, 9).Value = 0  # I59: 1000 -> 0
$ws.Cells.Item(59, 11).Value = 0  # K59: 1000 -> 0
$ws.Cells.Item(59, 13).ClearContents()  # M59: -417 -> (blank)

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Cells.Item(10, 8).Value = 673.6667  # H10: 875 -> 673.6667
$ws.Cells.Item(10, 9).Value = 408.4  # I10: 500 -> 408.4
$ws.Cells.Item(10, 10).Value = 2000  # J10: 1250 -> 2000
$ws.Cells.Item(10, 11).Value = 408.4  # K10: 500 -> 408.4
$ws.Cells.Item(10, 12).Value = 2000  # L10: 1250 -> 2000
$ws.Cells.Item(10, 13).Value = -268.4  # M10: -360 -> -268.4
$ws.Cells.Item(10, 14).Value = -2280  # N10: -1530 -> -2280
# Row 12
$ws.Cells.Item(12, 8).Value = 376.5  # H12: 403 -> 376.5
$ws.Cells.Item(12, 9).Value = 376.5  # I12: 403 -> 376.5
$ws.Cells.Item(12, 11).Value = 376.5  # K12: 403 -> 376.5
$ws.Cells.Item(12, 13).Value = -206.5  # M12: -233 -> -206.5
# Row 22
$ws.Cells.Item(22, 8).Value = 1333  # H22: 1337.25 -> 1333
$ws.Cells.Item(22, 10).Value = 1500  # J22: 1450 -> 1500
$ws.Cells.Item(22, 12).Value = 1500  # L22: 1450 -> 1500
$ws.Cells.Item(22, 14).Value = -2090  # N22: -2040 -> -2090
# Row 27
$ws.Cells.Item(27, 8).Value = 1333  # H27: 1337.25 -> 1333
$ws.Cells.Item(27, 10).Value = 1500  # J27: 1450 -> 1500
$ws.Cells.Item(27, 12).Value = 1500  # L27: 1450 -> 1500
$ws.Cells.Item(27, 14).Value = -1714  # N27: -1664 -> -1714
# Row 40
$ws.Cells.Item(40, 8).Value = 2328.0908  # H40: 2396.2 -> 2328.0908
$ws.Cells.Item(40, 9).Value = 2334.3333  # I40: 2420.25 -> 2334.3333
$ws.Cells.Item(40, 11).Value = 2334.3333  # K40: 2420.25 -> 2334.3333
$ws.Cells.Item(40, 13).Value = -2198.3333  # M40: -2284.25 -> -2198.3333
# Row 46
$ws.Cells.Item(46, 8).Value = 4416.6665  # H46: 1500 -> 4416.6665
$ws.Cells.Item(46, 10).Value = 5000  # J46: 0 -> 5000
$ws.Cells.Item(46, 12).Value = 5000  # L46: 0 -> 5000
$ws.Cells.Item(46, 14).Value = -5376  # N46: None -> -5376

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Cells.Item(3, 8).Value = 43809.875  # H3: 48354.855 -> 43809.875
$ws.Cells.Item(3, 9).Value = 54082.5  # I3: 62500 -> 54082.5
$ws.Cells.Item(3, 11).Value = 54082.5  # K3: 62500 -> 54082.5
$ws.Cells.Item(3, 13).Value = -53968.5  # M3: -62386 -> -53968.5
# Row 10
$ws.Cells.Item(10, 8).Value = 10000250  # H10: 10004000 -> 10000250
$ws.Cells.Item(10, 10).Value = 500  # J10: 8000 -> 500
$ws.Cells.Item(10, 12).Value = 500  # L10: 8000 -> 500
$ws.Cells.Item(10, 14).Value = -838  # N10: -8338 -> -838
# Row 12
$ws.Cells.Item(12, 8).Value = 39500  # H12: 0 -> 39500
$ws.Cells.Item(12, 9).Value = 39000  # I12: 0 -> 39000
$ws.Cells.Item(12, 10).Value = 40000  # J12: 0 -> 40000
$ws.Cells.Item(12, 11).Value = 39000  # K12: 0 -> 39000
$ws.Cells.Item(12, 12).Value = 40000  # L12: 0 -> 40000
$ws.Cells.Item(12, 13).Value = -38858  # M12: None -> -38858
$ws.Cells.Item(12, 14).Value = -40284  # N12: None -> -40284
# Row 54
$ws.Cells.Item(54, 8).Value = 24545.273  # H54: 24900 -> 24545.273
$ws.Cells.Item(54, 9).Value = 30000  # I54: 28166.666 -> 30000
$ws.Cells.Item(54, 10).Value = 19999.666  # J54: 20000 -> 19999.666
$ws.Cells.Item(54, 11).Value = 30000  # K54: 28166.666 -> 30000
$ws.Cells.Item(54, 12).Value = 19999.666  # L54: 20000 -> 19999.666
$ws.Cells.Item(54, 13).Value = -29480  # M54: -27646.666 -> -29480
$ws.Cells.Item(54, 14).Value = -21039.666  # N54: -21040 -> -21039.666
# Row 75
$ws.Cells.Item(75, 8).Value = 32500  # H75: 37559 -> 32500
$ws.Cells.Item(75, 9).Value = 0  # I75: 40118 -> 0
$ws.Cells.Item(75, 10).Value = 32500  # J75: 35000 -> 32500
$ws.Cells.Item(75, 11).Value = 0  # K75: 40118 -> 0
$ws.Cells.Item(75, 12).Value = 32500  # L75: 35000 -> 32500
$ws.Cells.Item(75, 13).ClearContents()  # M75: -39182 -> (blank)
$ws.Cells.Item(75, 14).Value = -34372  # N75: -36872 -> -34372
# Row 78
$ws.Cells.Item(78, 8).Value = 32500  # H78: 37559 -> 32500
$ws.Cells.Item(78, 9).Value = 0  # I78: 40118 -> 0
$ws.Cells.Item(78, 10).Value = 32500  # J78: 35000 -> 32500
$ws.Cells.Item(78, 11).Value = 97500  # K78: 120354 -> 97500
$ws.Cells.Item(78, 12).ClearContents()  # L78: 105000 -> (blank)
$ws.Cells.Item(78, 13).ClearContents()  # M78: -115674 -> (blank)
$ws.Cells.Item(78, 14).Value = -106860  # N78: -114360 -> -106860
# Row 132
$ws.Cells.Item(132, 8).Value = 1794.8889  # H132: 1965.875 -> 1794.8889
$ws.Cells.Item(132, 9).Value = 1744.25  # I132: 1932.4286 -> 1744.25
$ws.Cells.Item(132, 11).Value = 5232.75  # K132: 5797.2858 -> 5232.75
$ws.Cells.Item(132, 13).Value = -2702.75  # M132: -3267.2858 -> -2702.75
